$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), matching the existing
# header style (bold, centered, bordered) used by the other header cells.
$ws.Range("I1").Value = "I0"
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)

$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)

# Add new data cells I2 and J2 under the new headers.
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9
